$d = $word.ActiveDocument

# --- Part 1: insert the new content paragraphs before the bookmark paragraph ---
$bmRange = $d.Bookmarks.Item("_GoBack").Range
$bmPara = $bmRange.Paragraphs.Item(1)
$bmPara.Range.InsertParagraphBefore()

$bmRange2 = $d.Bookmarks.Item("_GoBack").Range
$bmPara2 = $bmRange2.Paragraphs.Item(1)
$idx2 = $bmPara2.Index
$placeholderPara = $d.Paragraphs.Item($idx2 - 1)
$xml1 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p/><w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">3. </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>索引列的顺序</w:t></w:r></w:p><w:p/><w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>让选择性最强的索引列放在前面。</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>索引的选择性是指：不重复的索引值和记录总数的比值。最大值为</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> 1</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>，此时每个记录都有唯一的索引与其对应。选择性越高，查询效率也越高。</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$placeholderPara.Range.InsertXML($xml1)

# --- Part 2: insert one more empty paragraph right after the bookmark paragraph ---
$bmRange3 = $d.Bookmarks.Item("_GoBack").Range
$bmPara3 = $bmRange3.Paragraphs.Item(1)
$idx3 = $bmPara3.Index
$nextPara = $d.Paragraphs.Item($idx3 + 1)
$xml2 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p/><w:p><w:pPr><w:jc w:val="left"/></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$nextPara.Range.InsertXML($xml2)

# Clean up the temporary alignment marker used to avoid an insertion no-op
$newEmptyPara = $d.Paragraphs.Item($idx3 + 2)
$newEmptyPara.Alignment = 0

Write-Output "DONE paragraphs=$($d.Paragraphs.Count)"
